# Update rows 2-10 (A:T) of Sheet1 with refreshed NATMI Rtn4-Rtn4rl1 values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Rtn4"
$ws.Cells.Item(2, 3).Value = "Rtn4rl1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 63.211268
$ws.Cells.Item(2, 8).Value = 189.633804
$ws.Cells.Item(2, 9).Value = 0.4922609885657722
$ws.Cells.Item(2, 10).Value = 0.4922609885657722
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.128483
$ws.Cells.Item(2, 14).Value = 0.385449
$ws.Cells.Item(2, 15).Value = 0.01580103135560779
$ws.Cells.Item(2, 16).Value = 0.01580103135560779
$ws.Cells.Item(2, 17).Value = 8.121573346443999
$ws.Cells.Item(2, 18).Value = 73.094160117996
$ws.Cells.Item(2, 19).Value = 0.007778231315470255
$ws.Cells.Item(2, 20).Value = 0.007778231315470255

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Rtn4"
$ws.Cells.Item(3, 3).Value = "Rtn4rl1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 63.211268
$ws.Cells.Item(3, 8).Value = 189.633804
$ws.Cells.Item(3, 9).Value = 0.4922609885657722
$ws.Cells.Item(3, 10).Value = 0.4922609885657722
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.550946666666667
$ws.Cells.Item(3, 14).Value = 13.65284
$ws.Cells.Item(3, 15).Value = 0.5596822223772701
$ws.Cells.Item(3, 16).Value = 0.55968222237727
$ws.Cells.Item(3, 17).Value = 287.6711094003734
$ws.Cells.Item(3, 18).Value = 2589.03998460336
$ws.Cells.Item(3, 19).Value = 0.2755097240701233
$ws.Cells.Item(3, 20).Value = 0.2755097240701233

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Rtn4"
$ws.Cells.Item(4, 3).Value = "Rtn4rl1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 63.211268
$ws.Cells.Item(4, 8).Value = 189.633804
$ws.Cells.Item(4, 9).Value = 0.4922609885657722
$ws.Cells.Item(4, 10).Value = 0.4922609885657722
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.451875
$ws.Cells.Item(4, 14).Value = 10.355625
$ws.Cells.Item(4, 15).Value = 0.4245167462671222
$ws.Cells.Item(4, 16).Value = 0.4245167462671222
$ws.Cells.Item(4, 17).Value = 218.1973957275
$ws.Cells.Item(4, 18).Value = 1963.7765615475
$ws.Cells.Item(4, 19).Value = 0.2089730331801786
$ws.Cells.Item(4, 20).Value = 0.2089730331801786

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Rtn4"
$ws.Cells.Item(5, 3).Value = "Rtn4rl1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 43.30706799999999
$ws.Cells.Item(5, 8).Value = 129.921204
$ws.Cells.Item(5, 9).Value = 0.3372560111523963
$ws.Cells.Item(5, 10).Value = 0.3372560111523963
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.128483
$ws.Cells.Item(5, 14).Value = 0.385449
$ws.Cells.Item(5, 15).Value = 0.01580103135560779
$ws.Cells.Item(5, 16).Value = 0.01580103135560779
$ws.Cells.Item(5, 17).Value = 5.564222017843998
$ws.Cells.Item(5, 18).Value = 50.077998160596
$ws.Cells.Item(5, 19).Value = 0.005328992807086226
$ws.Cells.Item(5, 20).Value = 0.005328992807086226

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Rtn4"
$ws.Cells.Item(6, 3).Value = "Rtn4rl1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 43.30706799999999
$ws.Cells.Item(6, 8).Value = 129.921204
$ws.Cells.Item(6, 9).Value = 0.3372560111523963
$ws.Cells.Item(6, 10).Value = 0.3372560111523963
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.550946666666667
$ws.Cells.Item(6, 14).Value = 13.65284
$ws.Cells.Item(6, 15).Value = 0.5596822223772701
$ws.Cells.Item(6, 16).Value = 0.55968222237727
$ws.Cells.Item(6, 17).Value = 197.0881567577067
$ws.Cells.Item(6, 18).Value = 1773.79341081936
$ws.Cells.Item(6, 19).Value = 0.1887561938318666
$ws.Cells.Item(6, 20).Value = 0.1887561938318665

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Rtn4"
$ws.Cells.Item(7, 3).Value = "Rtn4rl1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 43.30706799999999
$ws.Cells.Item(7, 8).Value = 129.921204
$ws.Cells.Item(7, 9).Value = 0.3372560111523963
$ws.Cells.Item(7, 10).Value = 0.3372560111523963
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.451875
$ws.Cells.Item(7, 14).Value = 10.355625
$ws.Cells.Item(7, 15).Value = 0.4245167462671222
$ws.Cells.Item(7, 16).Value = 0.4245167462671222
$ws.Cells.Item(7, 17).Value = 149.4905853525
$ws.Cells.Item(7, 18).Value = 1345.4152681725
$ws.Cells.Item(7, 19).Value = 0.1431708245134436
$ws.Cells.Item(7, 20).Value = 0.1431708245134436

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Rtn4"
$ws.Cells.Item(8, 3).Value = "Rtn4rl1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 21.891734
$ws.Cells.Item(8, 8).Value = 65.675202
$ws.Cells.Item(8, 9).Value = 0.1704830002818315
$ws.Cells.Item(8, 10).Value = 0.1704830002818315
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.128483
$ws.Cells.Item(8, 14).Value = 0.385449
$ws.Cells.Item(8, 15).Value = 0.01580103135560779
$ws.Cells.Item(8, 16).Value = 0.01580103135560779
$ws.Cells.Item(8, 17).Value = 2.812715659522
$ws.Cells.Item(8, 18).Value = 25.314440935698
$ws.Cells.Item(8, 19).Value = 0.002693807233051312
$ws.Cells.Item(8, 20).Value = 0.002693807233051311

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Rtn4"
$ws.Cells.Item(9, 3).Value = "Rtn4rl1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 21.891734
$ws.Cells.Item(9, 8).Value = 65.675202
$ws.Cells.Item(9, 9).Value = 0.1704830002818315
$ws.Cells.Item(9, 10).Value = 0.1704830002818315
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 4.550946666666667
$ws.Cells.Item(9, 14).Value = 13.65284
$ws.Cells.Item(9, 15).Value = 0.5596822223772701
$ws.Cells.Item(9, 16).Value = 0.55968222237727
$ws.Cells.Item(9, 17).Value = 99.62811387485334
$ws.Cells.Item(9, 18).Value = 896.65302487368
$ws.Cells.Item(9, 19).Value = 0.09541630447528021
$ws.Cells.Item(9, 20).Value = 0.09541630447528017

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Rtn4"
$ws.Cells.Item(10, 3).Value = "Rtn4rl1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 21.891734
$ws.Cells.Item(10, 8).Value = 65.675202
$ws.Cells.Item(10, 9).Value = 0.1704830002818315
$ws.Cells.Item(10, 10).Value = 0.1704830002818315
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 3.451875
$ws.Cells.Item(10, 14).Value = 10.355625
$ws.Cells.Item(10, 15).Value = 0.4245167462671222
$ws.Cells.Item(10, 16).Value = 0.4245167462671222
$ws.Cells.Item(10, 17).Value = 75.56752930124999
$ws.Cells.Item(10, 18).Value = 680.10776371125
$ws.Cells.Item(10, 19).Value = 0.07237288857349998
$ws.Cells.Item(10, 20).Value = 0.07237288857349997

